# rule_auto.xlsx — "mark red in template.xlsx"
#
# The authored diff boils down to three real edits on the "db_node" sheet
# (plus the resulting, automatic shared-string renumbering that ripples
# into every other sheet that referenced those strings):
#
#   1. Row 2 (NO.1 / msq_u_auto): the lone "manga" value that lived in I2
#      moves left into the previously-empty H2, and J2's duplicate
#      "manga" is removed — leaving H2="manga", K2="manga".
#   2. Rows 3 and 4 were orphan/leftover rows (192.168.55.252 / credential /
#      portfolio and 192.168.55.254) with no NO./db_name/node_name data —
#      delete them outright, shifting NO.2 (msq_c1_auto) and NO.3
#      (kfk_u_auto) up to rows 3 and 4.
#   3. The active sheet/tab moves from "db_node" to "msq_u" (the first
#      sheet), and the remembered selection on "db_node" becomes J27.
#
# Deleting the rows removes the now-unused shared strings
# ("192.168.55.252", "credential", "192.168.55.254") which is exactly what
# the sharedStrings.xml hunk shows (139 -> 136 unique strings); every other
# worksheet's <v> shared-string indices shift automatically to match.

$wb = $excel.ActiveWorkbook

$wsDbNode = $wb.Worksheets.Item("db_node")
$wsMsqU   = $wb.Worksheets.Item("msq_u")

# 1) Row 2: consolidate the stray "manga" cells (I2,J2,K2) -> (H2,K2)
$wsDbNode.Range("H2").Value = "manga"
$wsDbNode.Range("I2:J2").ClearContents()

# 2) Remove the two orphan rows (old rows 3 & 4); NO.2/NO.3 shift up
$wsDbNode.Rows("3:4").Delete()

# 3) Update the remembered selection / active tab:
#    - db_node keeps a plain selection at J27, no longer the active tab
#    - msq_u becomes the active tab (selection left at A1)
$wsDbNode.Range("J27").Select()
$wsMsqU.Activate()
$wsMsqU.Range("A1").Select()
